$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (last-changed) date for rows 2-3
$ws.Range("C2").Value2 = 46079
$ws.Range("C3").Value2 = 46079

# Rows 4-7: reorder the 4 "HÅBO" block records (A 47653-2024, A 37417-2023, A 4521-2024, A 35197-2025)
# and bump their "Förändrad" date.
# Row 4
$ws.Range("A4").Value2 = "A 47653-2024"
$ws.Range("B4").Value2 = 45588
$ws.Range("C4").Value2 = 46079
$ws.Range("G4").Value2 = 3
$ws.Range("H4").Value2 = 2
$ws.Range("I4").Value2 = 1
$ws.Range("J4").Value2 = 1
$ws.Range("O4").Value2 = 1
$ws.Range("R4").Value2 = "Backklöver`r`nNästrot`r`nBlåsippa"
$ws.Range("S4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/artfynd/A 47653-2024 artfynd.xlsx`", `"A 47653-2024`")"
$ws.Range("T4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/kartor/A 47653-2024 karta.png`", `"A 47653-2024`")"
$ws.Range("V4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/klagomål/A 47653-2024 FSC-klagomål.docx`", `"A 47653-2024`")"
$ws.Range("W4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/klagomålsmail/A 47653-2024 FSC-klagomål mail.docx`", `"A 47653-2024`")"
$ws.Range("X4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/tillsyn/A 47653-2024 tillsynsbegäran.docx`", `"A 47653-2024`")"
$ws.Range("Y4").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/tillsynsmail/A 47653-2024 tillsynsbegäran mail.docx`", `"A 47653-2024`")"
$ws.Range("Z4").ClearContents()

# Row 6
$ws.Range("A6").Value2 = "A 4521-2024"
$ws.Range("B6").Value2 = 45327
$ws.Range("C6").Value2 = 46079
$ws.Range("G6").Value2 = 17.4
$ws.Range("H6").Value2 = 1
$ws.Range("I6").Value2 = 2
$ws.Range("R6").Value2 = "Talltita`r`nBronshjon`r`nFjällig taggsvamp s.str."
$ws.Range("S6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/artfynd/A 4521-2024 artfynd.xlsx`", `"A 4521-2024`")"
$ws.Range("T6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/kartor/A 4521-2024 karta.png`", `"A 4521-2024`")"
$ws.Range("V6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/klagomål/A 4521-2024 FSC-klagomål.docx`", `"A 4521-2024`")"
$ws.Range("W6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/klagomålsmail/A 4521-2024 FSC-klagomål mail.docx`", `"A 4521-2024`")"
$ws.Range("X6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/tillsyn/A 4521-2024 tillsynsbegäran.docx`", `"A 4521-2024`")"
$ws.Range("Y6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/tillsynsmail/A 4521-2024 tillsynsbegäran mail.docx`", `"A 4521-2024`")"
$ws.Range("Z6").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/fåglar/A 4521-2024 prioriterade fågelarter.docx`", `"A 4521-2024`")"

# Row 7
$ws.Range("A7").Value2 = "A 35197-2025"
$ws.Range("B7").Value2 = 45853
$ws.Range("C7").Value2 = 46079
$ws.Range("G7").Value2 = 0.9
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 3
$ws.Range("O7").Value2 = 3
$ws.Range("R7").Value2 = "Grönsångare`r`nTallticka`r`nVintertagging"
$ws.Range("S7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/artfynd/A 35197-2025 artfynd.xlsx`", `"A 35197-2025`")"
$ws.Range("T7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/kartor/A 35197-2025 karta.png`", `"A 35197-2025`")"
$ws.Range("V7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/klagomål/A 35197-2025 FSC-klagomål.docx`", `"A 35197-2025`")"
$ws.Range("W7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/klagomålsmail/A 35197-2025 FSC-klagomål mail.docx`", `"A 35197-2025`")"
$ws.Range("X7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/tillsyn/A 35197-2025 tillsynsbegäran.docx`", `"A 35197-2025`")"
$ws.Range("Y7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/tillsynsmail/A 35197-2025 tillsynsbegäran mail.docx`", `"A 35197-2025`")"
$ws.Range("Z7").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_0305/fåglar/A 35197-2025 prioriterade fågelarter.docx`", `"A 35197-2025`")"

# Row 5 (A 37417-2023) keeps its position; only the "Förändrad" date changes
$ws.Range("C5").Value2 = 46079

# Rows 8-11: unaffected except "Förändrad" date bump
$ws.Range("C8").Value2 = 46079
$ws.Range("C9").Value2 = 46079
$ws.Range("C10").Value2 = 46079
$ws.Range("C11").Value2 = 46079

# Rows 12-31: permutation/refresh of the remaining records (A, B, "Förändrad", Markägare, Area)
# Row 12
$ws.Range("A12").Value2 = "A 37410-2023"
$ws.Range("B12").Value2 = 45155
$ws.Range("C12").Value2 = 46079
$ws.Range("G12").Value2 = 20.9

# Row 13
$ws.Range("A13").Value2 = "A 55962-2023"
$ws.Range("B13").Value2 = 45240
$ws.Range("C13").Value2 = 46079
$ws.Range("G13").Value2 = 3.4

# Row 14
$ws.Range("A14").Value2 = "A 12156-2023"
$ws.Range("B14").Value2 = 44998.49157407408
$ws.Range("C14").Value2 = 46079
$ws.Range("G14").Value2 = 0.5

# Row 15
$ws.Range("A15").Value2 = "A 12146-2023"
$ws.Range("B15").Value2 = 44998.47842592592
$ws.Range("C15").Value2 = 46079
$ws.Range("G15").Value2 = 3.1

# Row 16
$ws.Range("A16").Value2 = "A 32023-2023"
$ws.Range("B16").Value2 = 45119.49833333334
$ws.Range("C16").Value2 = 46079
$ws.Range("G16").Value2 = 3.1

# Row 17
$ws.Range("A17").Value2 = "A 11989-2025"
$ws.Range("B17").Value2 = 45728.60074074074
$ws.Range("C17").Value2 = 46079
$ws.Range("G17").Value2 = 9.6

# Row 18
$ws.Range("A18").Value2 = "A 23250-2022"
$ws.Range("B18").Value2 = 44719
$ws.Range("C18").Value2 = 46079
$ws.Range("F18").Value2 = "Naturvårdsverket"
$ws.Range("G18").Value2 = 1

# Row 19
$ws.Range("A19").Value2 = "A 15732-2025"
$ws.Range("B19").Value2 = 45747
$ws.Range("C19").Value2 = 46079
$ws.Range("F19").Value2 = "Kyrkan"
$ws.Range("G19").Value2 = 1.4

# Row 20
$ws.Range("A20").Value2 = "A 67005-2021"
$ws.Range("B20").Value2 = 44522
$ws.Range("C20").Value2 = 46079
$ws.Range("G20").Value2 = 1.3

# Row 21
$ws.Range("A21").Value2 = "A 23370-2025"
$ws.Range("B21").Value2 = 45791.70907407408
$ws.Range("C21").Value2 = 46079
$ws.Range("G21").Value2 = 3.8

# Row 22
$ws.Range("A22").Value2 = "A 45406-2025"
$ws.Range("B22").Value2 = 45922.42936342592
$ws.Range("C22").Value2 = 46079
$ws.Range("G22").Value2 = 9.1

# Row 23
$ws.Range("A23").Value2 = "A 34202-2022"
$ws.Range("B23").Value2 = 44791.64837962963
$ws.Range("C23").Value2 = 46079
$ws.Range("G23").Value2 = 2

# Row 24
$ws.Range("A24").Value2 = "A 4524-2024"
$ws.Range("B24").Value2 = 45327
$ws.Range("C24").Value2 = 46079
$ws.Range("G24").Value2 = 4.6

# Row 25
$ws.Range("A25").Value2 = "A 35300-2025"
$ws.Range("B25").Value2 = 45854.41511574074
$ws.Range("C25").Value2 = 46079
$ws.Range("G25").Value2 = 2.2

# Row 26
$ws.Range("A26").Value2 = "A 35198-2025"
$ws.Range("B26").Value2 = 45853
$ws.Range("C26").Value2 = 46079
$ws.Range("G26").Value2 = 1.2

# Row 27
$ws.Range("A27").Value2 = "A 35193-2025"
$ws.Range("B27").Value2 = 45853
$ws.Range("C27").Value2 = 46079
$ws.Range("G27").Value2 = 1.9

# Row 28
$ws.Range("A28").Value2 = "A 65836-2021"
$ws.Range("B28").Value2 = 44517
$ws.Range("C28").Value2 = 46079
$ws.Range("G28").Value2 = 1.8

# Row 29
$ws.Range("A29").Value2 = "A 12154-2023"
$ws.Range("B29").Value2 = 44998
$ws.Range("C29").Value2 = 46079
$ws.Range("G29").Value2 = 2.7

# Row 30
$ws.Range("A30").Value2 = "A 37415-2023"
$ws.Range("B30").Value2 = 45155
$ws.Range("C30").Value2 = 46079
$ws.Range("F30").ClearContents()
$ws.Range("G30").Value2 = 6.6

# Row 31
$ws.Range("A31").Value2 = "A 8486-2026"
$ws.Range("B31").Value2 = 46064
$ws.Range("C31").Value2 = 46079
$ws.Range("F31").ClearContents()
